# Generate and save output file after processing
#
# Inserts three new "general_college_subjects" columns (history, electives,
# cs) immediately before the existing "arts" column, shifting everything
# from the old column R onward three places to the right, then fills in
# the new header/data cells and normalizes a few values in row 2 that were
# re-cased ("Unknown"/"Considered" -> "unknown"/"considered") with H2 also
# changing from "Unknown" to "considered".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank columns at R:T, shifting general_college_subjects.arts
# (and everything after it) three columns to the right.
$ws.Range("R1:T1").EntireColumn.Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# New header cells for the inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data cells for row 2 (university of chicago) in the inserted columns.
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Normalize casing of existing row-2 values (and fix H2's stale "Unknown").
$ws.Range("D2").Value = "unknown"
$ws.Range("E2").Value = "considered"
$ws.Range("F2").Value = "considered"
$ws.Range("G2").Value = "considered"
$ws.Range("H2").Value = "considered"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"
